$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.387.49"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.875.44"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "238.74"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4797"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").Value = "0.06532"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "1.872.72"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "0.07490"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "16.55"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "5.065"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "88.40"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "0.6609"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "30.348.16"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "13.27"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "0.000007586"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "2.114.68"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").Value = "5.299"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "220.34"
$ws.Range("E23").Value = "  +14.65%  "
$ws.Range("D24").Value = "6.191"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "9.328"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "167.38"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "0.09395"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("D31").Value = "4.304"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").Value = "4.021"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "0.05015"
$ws.Range("D34").Value = "1.209"
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("D35").Value = "0.7422"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Value = "2.704"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D38").Value = "2.613"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "0.9057"
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "106.32"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "5.851"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "0.4270"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "7.424"
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("D46").Value = "64.40"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "0.1274"
$ws.Range("E47").Value = "  -7.89%  "
$ws.Range("D48").Value = "1.472"
$ws.Range("E48").Value = "  -7.15%  "
$ws.Range("D49").Value = "8.897"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").Value = "33.69"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").Value = "0.3880"
$ws.Range("E51").Value = "  +0.11%  "
